$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header block: new author + refreshed values ---
$ws.Range("B4").Value = "FRS"
$ws.Range("B5").Value = "Hemel"
$ws.Range("B6").Value = "DD-MM-YYYY"
$ws.Range("B7").Value = "DD-MM-YYYY"

# --- Fill in the new test scenario rows 26-31 (sheet rows 36-41). ---
# These rows already carry their target cell formatting in the template,
# so only the values need to be written.

$ws.Range("A36").Value = "TS_026"
$ws.Range("B36").Value = "FRS"
$ws.Range("C36").Value = "Validate the working of 'Newsletter' functionality"
$ws.Range("D36").Value = "P4"
$ws.Range("E36").Value = 13

$ws.Range("A37").Value = "TS_027"
$ws.Range("B37").Value = "FRS"
$ws.Range("C37").Value = "Validate the working of 'Contact Us' page functionality"
$ws.Range("D37").Value = "P4"
$ws.Range("E37").Value = 13

$ws.Range("A38").Value = "TS_028"
$ws.Range("B38").Value = "FRS"
$ws.Range("C38").Value = "Validate the working of 'Gift Certificate' page functionality"
$ws.Range("D38").Value = "P4"
$ws.Range("E38").Value = 11

$ws.Range("A39").Value = "TS_029"
$ws.Range("B39").Value = "FRS"
$ws.Range("C39").Value = "Validate the working of 'Speal Offers' page functionality"
$ws.Range("D39").Value = "P4"
$ws.Range("E39").Value = 16

$ws.Range("A40").Value = "TS_030"
$ws.Range("B40").Value = "FRS"
$ws.Range("C40").Value = "Validate the working of 'Header' options, 'Menu' options and 'Footer' options"
$ws.Range("D40").Value = "P4"
$ws.Range("E40").Value = 22

$ws.Range("A41").Value = "TS_031"
$ws.Range("B41").Value = "FRS"
$ws.Range("C41").Value = "Validate the complete Application functionality for different currencies"
$ws.Range("D41").Value = "P2"
$ws.Range("E41").Value = 3

# The "Number of Test Cases" column (E) in these new rows was authored with
# the same plain centred format as the adjacent Priority column (D), not the
# fill variant that the empty template cells carried - line up the format.
$ws.Range("D36").Copy()
$ws.Range("E36").PasteSpecial(-4122)
$ws.Range("D37").Copy()
$ws.Range("E37").PasteSpecial(-4122)
$ws.Range("D38").Copy()
$ws.Range("E38").PasteSpecial(-4122)
$ws.Range("D39").Copy()
$ws.Range("E39").PasteSpecial(-4122)
$ws.Range("D40").Copy()
$ws.Range("E40").PasteSpecial(-4122)
$ws.Range("D41").Copy()
$ws.Range("E41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New rows use the same 14pt row height as the rest of the data table.
$ws.Range("A36:E41").RowHeight = 14

# --- Sheet view: scroll position / zoom / selection left by the last save ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("B5:C5").Select()
